$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(302, 44376, 0, 0, "0")
    ,@(303, 44377, 0, 0, "0")
    ,@(304, 44378, 0, 0, "0")
    ,@(305, 44379, 0, 0, "0")
    ,@(306, 44380, 0, 0, "0")
    ,@(307, 44381, 0, 0, "0")
    ,@(308, 44382, 0, 0, "0")
    ,@(309, 44383, 0, 0, "0")
    ,@(310, 44384, 0, 0, "0")
    ,@(311, 44385, 0, 0, "0")
    ,@(312, 44386, 0, 0, "0")
    ,@(313, 44387, 0, 0, "0")
    ,@(314, 44388, 0, 0, "0")
    ,@(315, 44389, 0, 0, "0")
    ,@(316, 44390, 0, 0, "0")
    ,@(317, 44391, 0, 0, "0")
    ,@(318, 44392, 0, 0, "0")
    ,@(319, 44393, 0, 0, "0")
    ,@(320, 44394, 0, 0, "0")
    ,@(321, 44395, 0, 0, "0")
    ,@(322, 44396, 0, 0, "0")
    ,@(323, 44397, 0, 0, "0")
    ,@(324, 44398, 0, 0, "0")
    ,@(325, 44399, 1, 1, "11.49954001839926")
    ,@(326, 44400, 0, 1, "11.49954001839926")
    ,@(327, 44401, 2, 3, "34.4986200551978")
    ,@(328, 44402, 1, 4, "45.99816007359706")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = [double]$row[4]
}

# Match the date-column style (s="2": bold font, thin border, center/top
# alignment, custom datetime number format) used by existing column-A cells
# in rows 2-301, by copying the format from the last pre-existing row.
$ws.Range("A301").Copy() | Out-Null
$ws.Range("A302:A328").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
